$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efnb1", "Ephb1", "ECs", 3, 1, 7.484645, 22.453935, 0.3833353856819431, 0.3833353856819431, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 16.55503678733334, 148.995331086, 0.1844310076622395, 0.1844310076622395),
    @("ECs", "Efnb1", "Ephb1", "sCs", 3, 1, 7.484645, 22.453935, 0.3833353856819431, 0.3833353856819431, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 17.85420649714333, 160.68785847429, 0.1989043780197036, 0.1989043780197036),
    @("FAPs", "Efnb1", "Ephb1", "ECs", 3, 1, 8.299810000000001, 24.89943, 0.4250850731647057, 0.4250850731647057, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 18.35807307866667, 165.222657708, 0.2045176921156758, 0.2045176921156758),
    @("FAPs", "Efnb1", "Ephb1", "sCs", 3, 1, 8.299810000000001, 24.89943, 0.4250850731647057, 0.4250850731647057, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 19.79873749884667, 178.18863748962, 0.2205673810490298, 0.2205673810490299),
    @("M1", "Efnb1", "Ephb1", "ECs", 2, 0.6666666666666666, 0.493317, 1.479951, 0.02526584259620318, 0.02526584259620318, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 1.0911514284, 9.820362855600001, 0.01215594746402976, 0.01215594746402976),
    @("M1", "Efnb1", "Ephb1", "sCs", 2, 0.6666666666666666, 0.493317, 1.479951, 0.02526584259620318, 0.02526584259620318, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 1.176780406626, 10.591023659634, 0.01310989513217342, 0.01310989513217342),
    @("M2", "Efnb1", "Ephb1", "ECs", 2, 0.6666666666666666, 0.2771763333333333, 0.831529, 0.01419593001942512, 0.01419593001942512, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 0.6130770924888889, 5.5176938324, 0.006829971288790778, 0.006829971288790778),
    @("M2", "Efnb1", "Ephb1", "sCs", 2, 0.6666666666666666, 0.2771763333333333, 0.831529, 0.01419593001942512, 0.01419593001942512, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 0.6611888060762221, 5.950699254686, 0.007365958730634344, 0.007365958730634344),
    @("Neutro", "Efnb1", "Ephb1", "ECs", 3, 1, 0.9279860000000001, 2.783958, 0.04752795506232342, 0.04752795506232342, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 2.052581300533334, 18.4732317048, 0.02286673490545657, 0.02286673490545657),
    @("Neutro", "Efnb1", "Ephb1", "sCs", 3, 1, 0.9279860000000001, 2.783958, 0.04752795506232342, 0.04752795506232342, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 2.213659254441333, 19.922933289972, 0.02466122015686684, 0.02466122015686684),
    @("sCs", "Efnb1", "Ephb1", "ECs", 3, 1, 2.042122, 6.126366, 0.1045898134753994, 0.1045898134753994, 3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272, 4.516901581066668, 40.6521142296, 0.05032043847493473, 0.05032043847493474),
    @("sCs", "Efnb1", "Ephb1", "sCs", 3, 1, 2.042122, 6.126366, 0.1045898134753994, 0.1045898134753994, 3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727, 4.871369033582666, 43.842321302244, 0.05426937500046469, 0.0542693750004647)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
